$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 606, shifting existing rows 606:647 down to 607:648
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row 606 with the new entry
# Force column A to be treated as literal text (not auto-converted to a date serial),
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Cells.Item(606, 1).NumberFormat = "@"
$ws.Cells.Item(606, 1).Value = "2026/01/10"
$ws.Cells.Item(606, 1).Style = "Normal"
$ws.Cells.Item(606, 2).Value = "土"
$ws.Cells.Item(606, 3).Value = 15
$ws.Cells.Item(606, 4).Value = 201
